$d = $word.ActiveDocument

# Locate the (hidden) _GoBack bookmark that currently sits between the two
# runs of the "Explain that they will be adding..." bullet, splitting the
# word "day" into "d" / "ay".
$bm = $d.Bookmarks.Item("_GoBack")
$splitPos = $bm.Range.Start

# Split the paragraph into two paragraphs at that exact point. Word places
# the (collapsed) bookmark range at the start of the new, second paragraph.
$breakRange = $d.Range($splitPos, $splitPos)
$breakRange.InsertParagraphAfter()

# The first paragraph now ends with "...previous d" and the second begins
# with "ay by adding information about data." Re-join the tail onto the
# first paragraph so it reads as one complete sentence again.
$p1 = $d.Paragraphs.Item(32)
$p1EndRange = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$p1EndRange.InsertAfter("ay by adding information about data.")

# Replace the (now stale) leftover text in the second paragraph with the
# new accessibility note, keep the bookmark that already lives there, bump
# the paragraph to the next outline level, and highlight it green.
$p2 = $d.Paragraphs.Item(33)
$p2TextRange = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$newText = "We" + [char]0x2019 + "ve created an alternative list view of the chart in the " + [char]0x201C + "Revised" + [char]0x201D + " version of the chart file for visually impaired students."
$p2TextRange.Text = $newText

$p2 = $d.Paragraphs.Item(33)
$p2.Range.ListFormat.ListLevelNumber = 3
$p2.Range.HighlightColorIndex = 4
